$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17; this shifts existing rows 17:77 down to 18:78
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new record
$ws.Range("A17").Value2 = 2
$ws.Range("B17").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C17").Value2 = "Coquimbo"
$ws.Range("D17").Value2 = 44623
$ws.Range("E17").Value2 = 4
$ws.Range("F17").Value2 = 100112030
$ws.Range("G17").Value2 = "Poroto granado"
$ws.Range("H17").Value2 = "Sin especificar"
$ws.Range("I17").Value2 = "Primera"
$ws.Range("J17").Value2 = 900
$ws.Range("K17").Value2 = 23000
$ws.Range("L17").Value2 = 25000
$ws.Range("M17").Value2 = 24000
$ws.Range("N17").Value2 = "`$/malla 25 kilos"
$ws.Range("O17").Value2 = "Provincia de Limarí"
$ws.Range("P17").Value2 = 960
$ws.Range("Q17").Value2 = 25
$ws.Range("R17").Value2 = "Hortaliza"
